# Updated cryptos list on Fri May 31 21:07:35 UTC 2024 with GitHub Actions
# Column D = Price, Column E = Volume(1h)
# NOTE: many "Price" cells are stored as plain text that happens to look
# numeric (e.g. "1.00", "0.997"). Assigning such a string directly to
# .Value makes Excel auto-convert it to a number, which would change the
# cell's stored representation. A leading apostrophe forces Excel to keep
# the entered text as-is (exactly like typing '1.00 into a cell), so we
# use that for every Price value that would otherwise parse as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.671.33"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.792.01"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'594.95"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'166.53"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.792.03"
$ws.Range("E7").Value = "  +1.45%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.12%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.518"

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.23%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -1.95%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = "  +0.67%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  -0.60%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'36.27"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.431.42"
$ws.Range("E15").Value = "  +1.55%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.781.34"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "'18.58"
$ws.Range("E17").Value = "  +3.94%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "67.654.51"
$ws.Range("E18").Value = "  -1.12%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.35%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'6.99"
$ws.Range("E20").Value = "  -0.18%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'10.22"
$ws.Range("E21").Value = "  -4.18%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'457.36"
$ws.Range("E22").Value = "  -2.00%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.697"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - PEPE
$ws.Range("D24").Value = "'0.0000154"
$ws.Range("E24").Value = "  +7.77%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'83.46"
$ws.Range("E25").Value = "  -0.68%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'11.94"
$ws.Range("E26").Value = "  -0.91%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  -2.38%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'10.08"
$ws.Range("E28").Value = "  -0.20%  "

# Row 29 - Dai: unchanged

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.28%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.28"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +0.82%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'29.80"
$ws.Range("E33").Value = "  -0.19%  "

# Row 34 - Aptos
$ws.Range("D34").Value = "'9.19"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("E35").Value = "  -0.09%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "3.747.41"
$ws.Range("E36").Value = "  +1.43%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -1.17%  "

# Row 38 - dogwifhat
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  -2.17%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.27%  "

# Row 40 - Mantle
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - Filecoin
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42 - FirstDigitalUSD
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43 - USDe: unchanged

# Row 44 - Arweave
$ws.Range("D44").Value = "'44.91"
$ws.Range("E44").Value = "  +4.76%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -1.86%  "

# Row 46 - OKB
$ws.Range("D46").Value = "'47.12"
$ws.Range("E46").Value = "  +2.79%  "

# Row 47 - Cosmos
$ws.Range("D47").Value = "'8.38"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48 - Monero
$ws.Range("D48").Value = "'149.02"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49 & 50 - Bittensor and Stacks swap places (with updated data)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.84"
$ws.Range("E49").Value = "  -4.75%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'391.90"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51 - Maker
$ws.Range("D51").Value = "2.763.86"
$ws.Range("E51").Value = "  +2.56%  "
